# "Updates BRQSD -> CES"
# Turn on (1) additional sources as RPS-qualifying in the "RQSD-BRQSD" sheet
# to align with a Clean Energy Standard (CES) definition: nuclear, biomass
# w/ CCS, and small modular reactor all become qualifying sources across all
# forecast years (columns B:AE). The lignite row's shared "=B2" formulas get
# cleared to plain (still-zero) values.

$wb = $excel.ActiveWorkbook

$wsBRQSD = $wb.Worksheets.Item("RQSD-BRQSD")
$wsRQSD  = $wb.Worksheets.Item("RQSD-RQSD")

# nuclear (row 5): was all 0, now all 1 across B:AE
$wsBRQSD.Range("B5:AE5").Value = 1

# lignite (row 14): drop the shared "=B2" formulas, leaving the same
# (zero) literal values behind
$wsBRQSD.Range("B14:AE14").Value = 0

# biomass w CCS (row 21): was all 0, now all 1 across B:AE
$wsBRQSD.Range("B21:AE21").Value = 1

# small modular reactor (row 23): was all 0, now all 1 across B:AE
$wsBRQSD.Range("B23:AE23").Value = 1

# Selection / active-sheet bookkeeping: the workbook now opens focused on
# "RQSD-BRQSD" with E13 selected, while "RQSD-RQSD" keeps its prior
# selection but is no longer the tab in view.
$wsRQSD.Range("B24:AE25").Select()
$wsBRQSD.Activate()
$wsBRQSD.Range("E13").Select()
